$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 2, shifting existing row 2 (and its formatting) down to row 3
$ws.Rows.Item(2).Insert()

# Newly inserted row picks up formatting from the row above (header); clear it so it matches a plain data row
$ws.Range("A2:BD2").ClearFormats()

# 2. Rotate the Half-Time correct-score headers AW1:BC1 left by one position
#    (AW1 old value moves to BC1, others shift left) - BD1 (Odd_CS_4-4_HT) is unaffected
$ws.Range("AW1").Value = "Odd_CS_0-1_HT"
$ws.Range("AX1").Value = "Odd_CS_0-2_HT"
$ws.Range("AY1").Value = "Odd_CS_1-2_HT"
$ws.Range("AZ1").Value = "Odd_CS_0-3_HT"
$ws.Range("BA1").Value = "Odd_CS_1-3_HT"
$ws.Range("BB1").Value = "Odd_CS_2-3_HT"
$ws.Range("BC1").Value = "Odd_CS_3-3_HT"

# 3. Populate the new row 2 with the new match data (Bolivia - The Strongest vs GV San Jose)
$ws.Range("A2").Value = "MNIj0jM0"
$ws.Range("B2").Value = "30/10/2024"
$ws.Range("C2").Value = "10:00"
$ws.Range("D2").Value = "BOLIVIA - DIVISION PROFESIONAL"
$ws.Range("E2").Value = "The Strongest"
$ws.Range("F2").Value = "GV San Jose"
$ws.Range("G2").Value = 1.36
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 7.5
$ws.Range("J2").Value = 1.8
$ws.Range("K2").Value = 2.75
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 1.13
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 1.48
$ws.Range("R2").Value = 2.6
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 8
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 10
$ws.Range("AA2").Value = 11
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 10
$ws.Range("AE2").Value = 17
$ws.Range("AF2").Value = 41
$ws.Range("AG2").Value = 151
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 41
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 81
$ws.Range("AL2").Value = 51
$ws.Range("AM2").Value = 41
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 6.5
$ws.Range("AP2").Value = 15
$ws.Range("AQ2").Value = 17
$ws.Range("AR2").Value = 34
$ws.Range("AS2").Value = 81
$ws.Range("AT2").Value = 3.75
$ws.Range("AU2").Value = 8
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 8.5
$ws.Range("AX2").Value = 34
$ws.Range("AY2").Value = 34
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 101
$ws.Range("BB2").Value = 151
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51

# 4. Update row 3 (the match that was previously in row 2) for its few changed odds
#    and the left-rotation of its AW:BC half-time correct-score odds
$ws.Range("O3").Value = 1.18
$ws.Range("P3").Value = 4.5
$ws.Range("Q3").Value = 1.62
$ws.Range("R3").Value = 2.25
$ws.Range("AW3").Value = 5
$ws.Range("AX3").Value = 15
$ws.Range("AY3").Value = 19
$ws.Range("AZ3").Value = 41
$ws.Range("BA3").Value = 51
$ws.Range("BB3").Value = 101
$ws.Range("BC3").Value = 351
